# Append a new record row (row 5) to Sheet1, mirroring the structure of
# the existing data rows (2-4): base_id, date, type, number, amount,
# client_name, phone, location, note.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 5

$ws.Cells.Item($row, 1).Value = "20260101-002"

# Force the date column to stay literal text (like the existing rows,
# where "2025-12-08" etc. are plain strings, not date serials). Setting a
# text number format before the assignment keeps Excel from auto-parsing
# the string into a date; resetting the style afterwards drops the
# temporary formatting so the cell ends up with the same default style as
# its neighbours.
$dateCell = $ws.Cells.Item($row, 2)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026-01-01"
$dateCell.Style = "Normal"

$ws.Cells.Item($row, 3).Value = "q"
$ws.Cells.Item($row, 4).Value = "QUO-20260101-001"
$ws.Cells.Item($row, 5).Value = 0

# client_name / phone / note stay blank, like the rest of the sheet's
# empty cells (stored as empty text, not just a cleared/blank cell).
# Assigning a bare quote prefix ("'") forces Excel to store an empty
# text value instead of clearing the cell outright; resetting the style
# afterwards removes the transient quote-prefix formatting.
foreach ($col in 6, 7, 9) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'"
    $cell.Style = "Normal"
}

$ws.Cells.Item($row, 8).Value = "Abu Dhabi - Al Shamkha"
